$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing row (13) into the new row (14) so that the
# new row inherits the exact same cell styles (date format, centered text,
# etc.) used throughout the table.
$ws.Range("A13:E13").Copy($ws.Range("A14:E14"))

# Overwrite the two cells that differ for the new change-log entry.
$ws.Range("B14").Value = "12"
$ws.Range("D14").Value = "Initial responsibilities asigned in requirements document"

# Move the active selection down to A15, matching where the cursor would
# land after entering the new row of data.
$ws.Range("A15").Select()
